$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the existing table down by two rows (old row1 -> row3, etc.) by
#    inserting two blank rows at the very top.
# ---------------------------------------------------------------------------
$ws.Rows("1:2").Insert()

# ---------------------------------------------------------------------------
# 2. Build the four new fonts as a chain, so every intermediate state is
#    itself one of the fonts actually needed in the final workbook (avoids
#    leaving throw-away fonts in the palette):
#
#       default --(+Underline)--> U/Arial/10            "filler"  font
#               --(+Bold)-------> B+U/Arial/10           "label"   font
#               --(Name=Verdana)-> B+U/Verdana/10         "row2"   font
#               --(Size=14)------> B+U/Verdana/14         "title"  font
#
#    Each stage is produced on the cell that actually needs that exact
#    look, then copied (format only) to the other cells that share it.
# ---------------------------------------------------------------------------

# -- "filler" font: underline only, Arial 10 --------------------------------
$ws.Range("C19").Font.Underline = $true

# -- "label" font: bold + underline, Arial 10 --------------------------------
$ws.Range("C19").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Font.Bold = $true

# -- "row2" font: bold + underline, Verdana 10 -------------------------------
$ws.Range("A19").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Font.Name = "Verdana"

# -- "title" font: bold + underline, Verdana 14 ------------------------------
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Font.Size = 14

# ---------------------------------------------------------------------------
# 3. Propagate the four fonts (format-only copy) to every other cell that
#    needs them.
# ---------------------------------------------------------------------------
# filler -> rest of the row-19 banner line and the whole row-31 banner line
$ws.Range("C19").Copy()
$ws.Range("D19:J19").PasteSpecial(-4122)
$ws.Range("B31:J31").PasteSpecial(-4122)

# label -> A31
$ws.Range("A19").Copy()
$ws.Range("A31").PasteSpecial(-4122)

# title -> A17
$ws.Range("A1").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Values.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "macOS_Ver"
$ws.Range("A2").ClearContents()
$ws.Range("A17").Value = "Ver.__3.6.21.1._Uploaded_to_Mac-_App_Stpre_Connect"

$ws.Range("A19").Value = "TimeDateCalculator.macOS"
$ws.Range("B19").ClearContents()
$ws.Range("C19:J19").ClearContents()

$ws.Range("A31").Value = "TimeDateCalculatorDLL"
$ws.Range("B31:J31").ClearContents()

$ws.Range("A21").Value = "Xamarin.Essentials by Microsoft"
$ws.Range("G21").Value = "1.6.0-pre2"
$ws.Range("G22").Value = "1.6.1"
$ws.Range("A23").Value = "Xamarin.Essentials: a kit of essential API's for your apps"

$ws.Range("A26").Value = "Xamarin.Forms by Microsoft"
$ws.Range("G26").Value = "4.7.0.1351"
$ws.Range("G27").Value = "5.0.0.2012"
$ws.Range("A28").Value = "Build native Uis for iOS, Android, UWP, macOS, Tizen and many more from a single, shared C# codebase"

$ws.Range("A33").Value = "Xamarin.Essentials by Microsoft"
$ws.Range("G33").Value = "1.6.0-pre2"
$ws.Range("G34").Value = "1.6.1"
$ws.Range("A35").Value = "Xamarin.Essentials: a kit of essential API's for your apps"

$ws.Range("A38").Value = "Xamarin.Forms by Microsoft"
$ws.Range("G38").Value = "4.8.0.1560"
$ws.Range("G39").Value = "5.0.0.2012"
$ws.Range("A40").Value = "Build native UIs for iOS, Android, UWP, macOS, Tizen and many more from a single, shared C# codebase"

# ---------------------------------------------------------------------------
# 5. Row heights for the two big titles.
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 18
$ws.Rows(17).RowHeight = 18

# ---------------------------------------------------------------------------
# 6. Page setup + selection to match the saved view state.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("J13").Select()
